$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "95.635.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.80%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.573.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.37%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.53%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "652.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.89%  "

$ws.Range("E7").Value = "  +1.47%  "

$ws.Range("E8").Value = "  +0.69%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.01%  "

$ws.Range("E10").Value = "  -0.59%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.574.92"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.43%  "

$ws.Range("E12").Value = "  +0.32%  "

$ws.Range("E13").Value = "  -2.97%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.59"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.05%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.259.26"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.55%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "95.510.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.90%  "

$ws.Range("E17").Value = "  +0.56%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.575.38"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.71%  "

$ws.Range("E19").Value = "  -4.89%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.07%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.06%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.92%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "507.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.98%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.477"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.95%  "

$ws.Range("E25").Value = "  +3.56%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.62"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.72%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "95.53"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.85%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.58"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.50%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.763.08"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.49%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.96%  "

$ws.Range("E31").Value = "  -2.24%  "

$ws.Range("E32").Value = "  -0.05%  "

$ws.Range("E33").Value = "  -2.17%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.47%  "

$ws.Range("E35").Value = "  -1.56%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "32.01"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.63%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.562"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.45%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.18"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.43%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "563.91"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.61%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.50"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.42%  "

$ws.Range("E41").Value = "  -0.05%  "

$ws.Range("E42").Value = "  +0.37%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.904"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.82%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.78"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.71%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "35.10"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +37.51%  "

$ws.Range("E46").Value = "  +6.08%  "

$ws.Range("E47").Value = "  +1.98%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.59"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.74%  "

$ws.Range("E49").Value = "  -2.04%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.84%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.48"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.91%  "
